# Applies the commit "Update gh-pages to output generated at 456a3b4":
#   - refreshes the "想去人数" (want-to-go count, column F) figures that were
#     re-scraped on three of the four sheets, and
#   - appends one newly-scraped event ("北京·伦敦西区音乐剧明星演唱会（摇滚版）")
#     as row 28 of the "演出" (Performance) sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "演出" (Performance): append new row 28, then refresh column F
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Give the new index cell (A28) the same bold/centered/bordered style used
# by every other row-number cell in column A, by copying formats down from
# the row above (A27) before writing the real value.
$ws2.Range("A27").Copy() | Out-Null
$ws2.Range("A28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("A28").Value = 27

# B28 holds a literal text date ("2024-10-26"); format the cell as Text
# first so Excel's automatic date-recognition doesn't turn it into a
# datetime serial value, matching the original inline-string cell.
$ws2.Range("B28").NumberFormat = "@"
$ws2.Range("B28").Value = "2024-10-26"

$ws2.Range("C28").Value = "北京·伦敦西区音乐剧明星演唱会（摇滚版）"
$ws2.Range("D28").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws2.Range("E28").Value = "2024.10.26 14:30-10.26 16:30"
$ws2.Range("F28").Value = 0
$ws2.Range("G28").Value = 144
$ws2.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=89400"
$ws2.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202407/TYPRpfu21721116217467.jpeg"

# ------------------------------------------------------------------
# Sheet "展览" (Exhibition): refresh column F view/want-to-go counts
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8237
$ws1.Range("F6").Value = 107
$ws1.Range("F7").Value = 7196
$ws1.Range("F8").Value = 1138
$ws1.Range("F9").Value = 561
$ws1.Range("F11").Value = 718
$ws1.Range("F15").Value = 168
$ws1.Range("F17").Value = 100
$ws1.Range("F18").Value = 11813
$ws1.Range("F19").Value = 102
$ws1.Range("F20").Value = 8
$ws1.Range("F21").Value = 138
$ws1.Range("F22").Value = 2342
$ws1.Range("F24").Value = 3302
$ws1.Range("F27").Value = 2780
$ws1.Range("F28").Value = 107
$ws1.Range("F29").Value = 27
$ws1.Range("F31").Value = 3291
$ws1.Range("F33").Value = 2404
$ws1.Range("F35").Value = 1651
$ws1.Range("F37").Value = 109
$ws1.Range("F38").Value = 5888
$ws1.Range("F39").Value = 1250
$ws1.Range("F40").Value = 13
$ws1.Range("F41").Value = 165
$ws1.Range("F42").Value = 191
$ws1.Range("F43").Value = 1120
$ws1.Range("F44").Value = 1105
$ws1.Range("F45").Value = 1083
$ws1.Range("F46").Value = 1548
$ws1.Range("F47").Value = 11
$ws1.Range("F48").Value = 103

# ------------------------------------------------------------------
# Sheet "本地生活" (Local Life): refresh column F view/want-to-go counts
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 271
$ws3.Range("F3").Value = 415

# ------------------------------------------------------------------
# Sheet "全部类型" (All Types): refresh column F view/want-to-go counts
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 271
$ws4.Range("F5").Value = 415
$ws4.Range("F8").Value = 8237
$ws4.Range("F10").Value = 107
$ws4.Range("F11").Value = 7196
$ws4.Range("F12").Value = 7196
$ws4.Range("F13").Value = 1138
$ws4.Range("F14").Value = 561
$ws4.Range("F16").Value = 718
$ws4.Range("F20").Value = 100
$ws4.Range("F22").Value = 11814
$ws4.Range("F23").Value = 102
$ws4.Range("F25").Value = 138
$ws4.Range("F26").Value = 2342
$ws4.Range("F27").Value = 2342
$ws4.Range("F28").Value = 3302
$ws4.Range("F29").Value = 2780
$ws4.Range("F30").Value = 107
$ws4.Range("F31").Value = 27
$ws4.Range("F33").Value = 3291
$ws4.Range("F36").Value = 2404
$ws4.Range("F38").Value = 1651
$ws4.Range("F39").Value = 109
$ws4.Range("F40").Value = 5888
$ws4.Range("F43").Value = 1250
$ws4.Range("F44").Value = 165
$ws4.Range("F45").Value = 191
$ws4.Range("F46").Value = 1120
$ws4.Range("F47").Value = 1105
$ws4.Range("F48").Value = 1083
$ws4.Range("F49").Value = 1548
$ws4.Range("F50").Value = 103
